$d = $word.ActiveDocument

# 1 & 2: Replace scenario title text (occurs twice: TOC entry and heading)
$d.Content.Find.Execute("brisanje reportovanih pitanja", $true, $false, $false, $false, $false,
                         $true, 1, $false, "prikaz liste najboljih rezultata", 2)

# 3: Replace the "kratak opis" paragraph content
$d.Content.Find.Execute("Administrator briše neka od repostovanih pitanja. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Plyer pregleda listu najboljih rezultata. Rezultate ", 2)

# 4: Replace "Reports" with "High Score List"
$d.Content.Find.Execute("Reports", $true, $false, $false, $false, $false,
                         $true, 1, $false, "High Score List", 2)

# 5: Replace "Prikazuje se lista svih reportovanih pitanja"
$d.Content.Find.Execute("Prikazuje se lista svih reportovanih pitanja", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Prikazuje se lista svih Playera i njihovih rezultata", 2)

# 6: Replace "Brisanje pitanja pritiskom na 🗑" with new sort description
$d.Content.Find.Execute("Brisanje pitanja pritiskom na ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Pritiskom na zaglavlje neke kolone sortira se po njoj ulazno/silazno.", 2)

# 7: Replace "admin" with "player"
$d.Content.Find.Execute("Korisnik je ulogovan u sistem kao admin. ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Korisnik je ulogovan u sistem kao player. ", 2)

# 8: Replace "Pitanje je izbrisano iz baze." with "Nema."
$d.Content.Find.Execute("Pitanje je izbrisano iz baze.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Nema.", 2)
